$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: 'Bitcoin' -> 'Bitcoin'
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "59.452.73"
$ws.Range("E2").Value = "  +0.23%  "

# Row 3: 'Ethereum' -> 'Ethereum'
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.645.72"
$ws.Range("E3").Value = "  -0.47%  "

# Row 4: 'TetherUSD' -> 'TetherUSD'
$ws.Range("E4").Value = "  +0.16%  "

# Row 5: 'BNB' -> 'BNB'
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "531.20"
$ws.Range("E5").Value = "  +1.05%  "

# Row 6: 'Solana' -> 'Solana'
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "146.53"
$ws.Range("E6").Value = "  -0.32%  "

# Row 7: 'USDC' -> 'USDC'
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.999"
$ws.Range("E7").Value = "  +0.00%  "

# Row 8: 'XRP' -> 'XRP'
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.571"
$ws.Range("E8").Value = "  -0.79%  "

# Row 9: 'Toncoin' -> 'Toncoin'
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "6.69"
$ws.Range("E9").Value = "  -4.18%  "

# Row 10: 'Dogecoin' -> 'Dogecoin'
$ws.Range("E10").Value = "  +0.86%  "

# Row 11: 'Cardano' -> 'Cardano'
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.338"
$ws.Range("E11").Value = "  +0.07%  "

# Row 12: 'TRON' -> 'TRON'
$ws.Range("E12").Value = "  +0.40%  "

# Row 13: 'WrappedliquidstakedEther2.0' -> 'WrappedliquidstakedEther2.0'
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "3.132.62"
$ws.Range("E13").Value = "  +0.17%  "

# Row 14: 'WrappedBTC' -> 'WrappedBTC'
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "59.469.18"
$ws.Range("E14").Value = "  +0.24%  "

# Row 15: 'Avalanche' -> 'Avalanche'
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "20.90"
$ws.Range("E15").Value = "  -2.16%  "

# Row 16: 'ShibaInu' -> 'ShibaInu'
$ws.Range("E16").Value = "  -0.09%  "

# Row 17: 'WrappedEther' -> 'WrappedEther'
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "2.658.53"
$ws.Range("E17").Value = "  -0.47%  "

# Row 18: 'BitcoinCash' -> 'BitcoinCash'
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "343.03"
$ws.Range("E18").Value = "  +0.44%  "

# Row 19: 'Polkadot' -> 'Polkadot'
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "4.45"
$ws.Range("E19").Value = "  +0.55%  "

# Row 20: 'Chainlink' -> 'Chainlink'
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "10.65"
$ws.Range("E20").Value = "  +1.94%  "

# Row 21: 'Uniswap' -> 'Uniswap'
$ws.Range("E21").Value = "  +0.51%  "

# Row 22: 'Dai' -> 'Dai'
$ws.Range("E22").Value = "  -0.08%  "

# Row 23: 'Litecoin' -> 'Litecoin'
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "65.76"
$ws.Range("E23").Value = "  +2.73%  "

# Row 24: 'Polygon' -> 'Polygon'
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.417"
$ws.Range("E24").Value = "  +0.66%  "

# Row 25: 'Kaspa' -> 'Kaspa'
$ws.Range("E25").Value = "  +0.34%  "

# Row 26: 'Binance-PegBSC-USD' -> 'WrappedeETH'
$ws.Range("B26").Value = "WrappedeETH"
$ws.Range("C26").Value = "https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.767.35"
$ws.Range("E26").Value = "  -0.42%  "

# Row 27: 'InternetComputer(DFINITY)' -> 'Binance-PegBSC-USD'
$ws.Range("B27").Value = "Binance-PegBSC-USD"
$ws.Range("C27").Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.997"
$ws.Range("E27").Value = "  -0.16%  "

# Row 28: 'PEPE' -> 'InternetComputer(DFINITY)'
$ws.Range("B28").Value = "InternetComputer(DFINITY)"
$ws.Range("C28").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "7.20"
$ws.Range("E28").Value = "  +0.22%  "

# Row 29: 'USDe' -> 'PEPE'
$ws.Range("B29").Value = "PEPE"
$ws.Range("C29").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.0₃0802"
$ws.Range("E29").Value = "  -1.07%  "

# Row 30: 'Aptos' -> 'USDe'
$ws.Range("B30").Value = "USDe"
$ws.Range("C30").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.999"
$ws.Range("E30").Value = "  +0.06%  "

# Row 31: 'PancakeSwap' -> 'Aptos'
$ws.Range("B31").Value = "Aptos"
$ws.Range("C31").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "6.39"
$ws.Range("E31").Value = "  -3.42%  "

# Row 32: 'EthereumClassic' -> 'PancakeSwap'
$ws.Range("B32").Value = "PancakeSwap"
$ws.Range("C32").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.62"
$ws.Range("E32").Value = "  +1.13%  "

# Row 33: 'Monero' -> 'EthereumClassic'
$ws.Range("B33").Value = "EthereumClassic"
$ws.Range("C33").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "19.04"
$ws.Range("E33").Value = "  +1.05%  "

# Row 34: 'NEARProtocol' -> 'Monero'
$ws.Range("B34").Value = "Monero"
$ws.Range("C34").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "151.73"
$ws.Range("E34").Value = "  +1.72%  "

# Row 35: 'ImmutableX' -> 'NEARProtocol'
$ws.Range("B35").Value = "NEARProtocol"
$ws.Range("C35").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "4.19"
$ws.Range("E35").Value = "  -0.48%  "

# Row 36: 'SuiNetwork' -> 'ImmutableX'
$ws.Range("B36").Value = "ImmutableX"
$ws.Range("C36").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.18"
$ws.Range("E36").Value = "  -1.43%  "

# Row 37: 'Fetch.AI' -> 'SuiNetwork'
$ws.Range("B37").Value = "SuiNetwork"
$ws.Range("C37").Value = "https://coinranking.com/coin/3xJluUMvp+suinetwork-sui"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.859"
$ws.Range("E37").Value = "  -4.67%  "

# Row 38: 'Stacks' -> 'Fetch.AI'
$ws.Range("B38").Value = "Fetch.AI"
$ws.Range("C38").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.857"
$ws.Range("E38").Value = "  -4.36%  "

# Row 39: 'OKB' -> 'Stacks'
$ws.Range("B39").Value = "Stacks"
$ws.Range("C39").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.48"
$ws.Range("E39").Value = "  -0.53%  "

# Row 40: 'Filecoin' -> 'OKB'
$ws.Range("B40").Value = "OKB"
$ws.Range("C40").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "36.53"
$ws.Range("E40").Value = "  -0.70%  "

# Row 41: 'FirstDigitalUSD' -> 'Filecoin'
$ws.Range("B41").Value = "Filecoin"
$ws.Range("C41").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "3.64"
$ws.Range("E41").Value = "  +0.42%  "

# Row 42: 'Stellar' -> 'FirstDigitalUSD'
$ws.Range("B42").Value = "FirstDigitalUSD"
$ws.Range("C42").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.997"
$ws.Range("E42").Value = "  -0.14%  "

# Row 43: 'Mantle' -> 'Stellar'
$ws.Range("B43").Value = "Stellar"
$ws.Range("C43").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.0982"
$ws.Range("E43").Value = "  +0.44%  "

# Row 44: 'Bittensor' -> 'Mantle'
$ws.Range("B44").Value = "Mantle"
$ws.Range("C44").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.602"
$ws.Range("E44").Value = "  -3.25%  "

# Row 45: 'EnergySwap' -> 'Bittensor'
$ws.Range("B45").Value = "Bittensor"
$ws.Range("C45").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "271.62"
$ws.Range("E45").Value = "  -1.67%  "

# Row 46: 'WhiteBITCoin' -> 'EnergySwap'
$ws.Range("B46").Value = "EnergySwap"
$ws.Range("C46").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "19.43"
$ws.Range("E46").Value = "  -2.47%  "

# Row 47: 'Hedera' -> 'WhiteBITCoin'
$ws.Range("B47").Value = "WhiteBITCoin"
$ws.Range("C47").Value = "https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "10.71"
$ws.Range("E47").Value = "  +1.71%  "

# Row 48: 'Maker' -> 'Hedera'
$ws.Range("B48").Value = "Hedera"
$ws.Range("C48").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.0535"
$ws.Range("E48").Value = "  -1.34%  "

# Row 49: 'RenderToken' -> 'Maker'
$ws.Range("B49").Value = "Maker"
$ws.Range("C49").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.041.56"
$ws.Range("E49").Value = "  -0.72%  "

# Row 50: 'VeChain' -> 'RenderToken'
$ws.Range("B50").Value = "RenderToken"
$ws.Range("C50").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "4.77"
$ws.Range("E50").Value = "  -1.60%  "

# Row 51: 'InjectiveProtocol' -> 'VeChain'
$ws.Range("B51").Value = "VeChain"
$ws.Range("C51").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0229"
$ws.Range("E51").Value = "  -0.32%  "
